# fixed id vs rank bug
#
# Two Division-1-Women teams (British Columbia / Thunderbirds and
# Carleton College / Syzygy) were missing from the seed sheet entirely -
# they'd been dropped when the sheet's row "id" got confused with its
# finish "rank" during a prior edit. Re-append them as rows 22-23 with
# their finish/rank values, matching the existing sheet conventions:
#   - every column is plain text (the "finish"/"rank" columns hold
#     numeric-looking values like "1", "2" but are stored as text
#     throughout the sheet, e.g. E2/F2 on the existing rows), so the two
#     numeric-looking columns are entered with a leading apostrophe to
#     force text entry instead of Excel's automatic number inference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "British Columbia"
$ws.Range("B22").Value = "Thunderbirds"
$ws.Range("C22").Value = "British Columbia,UBC,TBirds,Thunderbirds"
$ws.Range("D22").Value = "d1w"
$ws.Range("E22").Value = "'1"
$ws.Range("F22").Value = "'1"

$ws.Range("A23").Value = "Carleton College"
$ws.Range("B23").Value = "Syzygy"
$ws.Range("C23").Value = "Carleton,Syzygy"
$ws.Range("D23").Value = "d1w"
$ws.Range("E23").Value = "'2"
$ws.Range("F23").Value = "'2"
